$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "'69.324.70"
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -2.72%  '
$ws.Cells.Item(3, 4).Value = "'3.684.37"
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -3.45%  '
$ws.Cells.Item(4, 5).Value = '  +0.07%  '
$ws.Cells.Item(5, 4).Value = "'688.52"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(6, 4).Value = "'162.48"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -5.69%  '
$ws.Cells.Item(7, 4).Value = "'3.682.54"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -3.48%  '
$ws.Cells.Item(8, 5).Value = '  +0.05%  '
$ws.Cells.Item(9, 5).Value = '  -4.62%  '
$ws.Cells.Item(10, 5).Value = '  -8.69%  '
$ws.Cells.Item(11, 4).Value = "'7.39"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -3.76%  '
$ws.Cells.Item(12, 5).Value = '  -3.88%  '
$ws.Cells.Item(13, 5).Value = '  -5.39%  '
$ws.Cells.Item(14, 4).Value = "'33.51"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -6.97%  '
$ws.Cells.Item(15, 4).Value = "'4.305.48"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -3.47%  '
$ws.Cells.Item(16, 4).Value = "'3.679.79"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -3.56%  '
$ws.Cells.Item(17, 4).Value = "'69.372.35"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -2.59%  '
$ws.Cells.Item(18, 5).Value = '  -1.09%  '
$ws.Cells.Item(19, 4).Value = "'16.30"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -7.17%  '
$ws.Cells.Item(20, 4).Value = "'6.61"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -7.86%  '
$ws.Cells.Item(21, 4).Value = "'481.77"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -6.05%  '
$ws.Cells.Item(22, 4).Value = "'9.94"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -7.17%  '
$ws.Cells.Item(23, 5).Value = '  -7.86%  '
$ws.Cells.Item(24, 4).Value = "'80.11"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -4.90%  '
$ws.Cells.Item(25, 4).Value = "'3.830.53"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -3.43%  '
$ws.Cells.Item(26, 4).Value = "'0.0000130"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -9.90%  '
$ws.Cells.Item(27, 5).Value = '  +0.00%  '
$ws.Cells.Item(28, 4).Value = "'11.44"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -5.03%  '
$ws.Cells.Item(29, 4).Value = "'9.51"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -8.79%  '
$ws.Cells.Item(30, 5).Value = '  -10.33%  '
$ws.Cells.Item(31, 5).Value = '  -10.20%  '
$ws.Cells.Item(32, 5).Value = '  -7.90%  '
$ws.Cells.Item(33, 5).Value = '  -7.58%  '
$ws.Cells.Item(34, 4).Value = "'27.10"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -7.04%  '
$ws.Cells.Item(36, 5).Value = '  -3.93%  '
$ws.Cells.Item(37, 4).Value = "'3.650.32"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -3.37%  '
$ws.Cells.Item(38, 4).Value = "'8.51"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -7.43%  '
$ws.Cells.Item(39, 4).Value = "'6.32"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +5.11%  '
$ws.Cells.Item(40, 5).Value = '  -2.16%  '
$ws.Cells.Item(41, 5).Value = '  -7.86%  '
$ws.Cells.Item(42, 5).Value = '  -0.02%  '
$ws.Cells.Item(43, 5).Value = '  +0.04%  '
$ws.Cells.Item(44, 4).Value = "'0.951"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -7.23%  '
$ws.Cells.Item(45, 4).Value = "'163.49"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -4.10%  '
$ws.Cells.Item(46, 4).Value = "'47.98"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -3.30%  '
$ws.Cells.Item(47, 4).Value = "'2.83"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -14.82%  '
$ws.Cells.Item(48, 4).Value = "'29.90"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +2.33%  '
$ws.Cells.Item(49, 4).Value = "'1.36"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +0.28%  '
$ws.Cells.Item(50, 4).Value = "'0.000287"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -7.71%  '
$ws.Cells.Item(51, 5).Value = '  -0.59%  '
